$wb = $excel.ActiveWorkbook

# ===== Sheet 1: 土地 (land) =====
$ws1 = $wb.Worksheets.Item(1)

# Extend formatting from column H into the new columns I:O
$ws1.Range("H1").Copy() | Out-Null
$ws1.Range("I1:O1").PasteSpecial(-4122) | Out-Null
$ws1.Range("H2").Copy() | Out-Null
$ws1.Range("I2:O2").PasteSpecial(-4122) | Out-Null
$ws1.Range("H3").Copy() | Out-Null
$ws1.Range("I3:O3").PasteSpecial(-4122) | Out-Null
$ws1.Range("H4").Copy() | Out-Null
$ws1.Range("I4:O4").PasteSpecial(-4122) | Out-Null

# Header values
$ws1.Range("B1").Value = "name"
$ws1.Range("C1").Value = "area"
$ws1.Range("D1").Value = "share_portion"
$ws1.Range("E1").Value = "owner"
$ws1.Range("F1").Value = "register_date"
$ws1.Range("G1").Value = "register_reason"
$ws1.Range("H1").Value = "acquire_value"
$ws1.Range("I1").Value = "property_category"
$ws1.Range("J1").Value = "category"
$ws1.Range("K1").Value = "date"
$ws1.Range("L1").Value = "legislator_name"
$ws1.Range("M1").Value = "legislator_id"
$ws1.Range("N1").Value = "source_file"
$ws1.Range("O1").Value = "index"

# Row 2
$ws1.Range("A2").Value = 15
$ws1.Range("B2").Value = "南投縣南投市牛運堀段02670002地號"
$ws1.Range("C2").Value = 325
$ws1.Range("D2").Value = "全部"
$ws1.Range("E2").Value = "廖述嘉"
$ws1.Range("F2").Value = "79年05月04日"
$ws1.Range("G2").Value = "共有物分割"
$ws1.Range("H2").Value = "(超過五年）"
$ws1.Range("I2").Value = "land"
$ws1.Range("J2").Value = "normal"
$ws1.Range("K2").Value = "2011-11-21"
$ws1.Range("L2").Value = "盧秀燕"
$ws1.Range("M2").Value = 869
$ws1.Range("N2").Value = "tmp9eb41"
$ws1.Range("O2").Value = 15

# Row 3
$ws1.Range("A3").Value = 16
$ws1.Range("B3").Value = "南投縣南投市牛運堀段02670008地號"
$ws1.Range("C3").Value = 27
$ws1.Range("D3").Value = "全部"
$ws1.Range("E3").Value = "廖述嘉"
$ws1.Range("F3").Value = "80年06月25日"
$ws1.Range("G3").Value = "共有物分割"
$ws1.Range("H3").Value = "(超過五年）"
$ws1.Range("I3").Value = "land"
$ws1.Range("J3").Value = "normal"
$ws1.Range("K3").Value = "2011-11-21"
$ws1.Range("L3").Value = "盧秀燕"
$ws1.Range("M3").Value = 869
$ws1.Range("N3").Value = "tmp9eb41"
$ws1.Range("O3").Value = 16

# Row 4
$ws1.Range("A4").Value = 17
$ws1.Range("B4").Value = "南投縣南投市牛運堀段02670041地號"
$ws1.Range("C4").Value = "3"
$ws1.Range("D4").Value = "全部"
$ws1.Range("E4").Value = "廖述嘉"
$ws1.Range("F4").Value = "80年06月25日"
$ws1.Range("G4").Value = "共有物分‘割"
$ws1.Range("H4").Value = "(超過五年）"
$ws1.Range("I4").Value = "land"
$ws1.Range("J4").Value = "normal"
$ws1.Range("K4").Value = "2011-11-21"
$ws1.Range("L4").Value = "盧秀燕"
$ws1.Range("M4").Value = 869
$ws1.Range("N4").Value = "tmp9eb41"
$ws1.Range("O4").Value = 17

# ===== Sheet 2: 汽車 (car) =====
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B1").Value = "廠牌型號"
$ws2.Range("C1").Value = "汽缸容量"
$ws2.Range("D1").Value = "所有人"
$ws2.Range("E1").Value = "登記（取得)時間"
$ws2.Range("F1").Value = "登記（取得）原因"
$ws2.Range("G1").Value = "取得價額"

$ws2.Range("B2").Value = "HYUNDAI"
$ws2.Range("C2").Value = 2497
$ws2.Range("D2").Value = "盧秀燕"
$ws2.Range("E2").Value = "99年02月06日"
$ws2.Range("F2").Value = "(購二手車）"
$ws2.Range("G2").Value = 100000

$ws2.Range("B3").Value = "國瑞"
$ws2.Range("C3").Value = 1998
$ws2.Range("D3").Value = "廖述嘉"
$ws2.Range("E3").Value = "99年08月16日"
$ws2.Range("F3").Value = "(購二手車）"
$ws2.Range("G3").Value = 50000

# ===== Sheet 3: 存款 (deposit) =====
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("B1").Value = "存放機構(應敘明分支機構）"
$ws3.Range("C1").Value = "種類"
$ws3.Range("D1").Value = "幣別"
$ws3.Range("E1").Value = "所有人"
$ws3.Range("F1").Value = "新臺幣總額或折合新臺幣總額"

$ws3.Range("A2").Value = 47
$ws3.Range("B2").Value = "臺灣銀行群賢分行"
$ws3.Range("C2").Value = "活期儲蓄存款"
$ws3.Range("D2").Value = "新臺幣"
$ws3.Range("E2").Value = "盧秀燕"
$ws3.Range("F2").Value = 4752062

$ws3.Range("A3").Value = 48
$ws3.Range("B3").Value = "臺灣銀行群賢分行"
$ws3.Range("C3").Value = "定期存款"
$ws3.Range("D3").Value = "新臺幣"
$ws3.Range("E3").Value = "盧秀燕"
$ws3.Range("F3").Value = 2000000

$ws3.Range("A4").Value = 49
$ws3.Range("B4").Value = "臺灣銀行群賢分行"
$ws3.Range("C4").Value = "活期儲蓄存款"
$ws3.Range("D4").Value = "新臺幣"
$ws3.Range("E4").Value = "盧秀燕"
$ws3.Range("F4").Value = 188252

$ws3.Range("A5").Value = 50
$ws3.Range("B5").Value = "臺灣銀行群賢分行"
$ws3.Range("C5").Value = "支票存款"
$ws3.Range("D5").Value = "新臺幣"
$ws3.Range("E5").Value = "盧秀燕"
$ws3.Range("F5").Value = 37158

$ws3.Range("A6").Value = 51
$ws3.Range("B6").Value = "合作金庫商業銀行中權分行"
$ws3.Range("C6").Value = "活期儲蓄存款"
$ws3.Range("D6").Value = "新臺幣"
$ws3.Range("E6").Value = "廖〇青"
$ws3.Range("F6").Value = 155001

$ws3.Range("A7").Value = 52
$ws3.Range("B7").Value = "中國信託商業銀行城中分行"
$ws3.Range("C7").Value = "活期儲蓄存款"
$ws3.Range("D7").Value = "新臺幣"
$ws3.Range("E7").Value = "廖〇青"
$ws3.Range("F7").Value = 3531

$ws3.Range("A8").Value = 53
$ws3.Range("B8").Value = "中華郵政股份有限公司中正路郵局"
$ws3.Range("C8").Value = "活期存款"
$ws3.Range("D8").Value = "新臺幣"
$ws3.Range("E8").Value = "廖〇青"
$ws3.Range("F8").Value = 805752

$ws3.Range("A9").Value = 54
$ws3.Range("B9").Value = "中華郵政股份有限公司中正路郵局"
$ws3.Range("C9").Value = "定期存款"
$ws3.Range("D9").Value = "新臺幣"
$ws3.Range("E9").Value = "廖〇青"
$ws3.Range("F9").Value = 840000

$ws3.Range("A10").Value = 55
$ws3.Range("B10").Value = "中華郵政股份有限公司台中永安郵局"
$ws3.Range("C10").Value = "活期存款"
$ws3.Range("D10").Value = "新臺幣"
$ws3.Range("E10").Value = "廖〇桐"
$ws3.Range("F10").Value = 1238149

$ws3.Range("A11").Value = 56
$ws3.Range("B11").Value = "國泰世華商業銀行中港分行"
$ws3.Range("C11").Value = "活期儲蓄存款"
$ws3.Range("D11").Value = "新臺幣"
$ws3.Range("E11").Value = "廖〇桐"
$ws3.Range("F11").Value = 45421

# ===== Sheet 4: 保險 (insurance) =====
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("B1").Value = "保險公司"
$ws4.Range("C1").Value = "保險名稱"
$ws4.Range("D1").Value = "要保人"
$ws4.Range("E1").Value = "備註"

$ws4.Range("A2").Value = 87
$ws4.Range("B2").Value = "南山人壽"
$ws4.Range("C2").Value = "子女教育保險"
$ws4.Range("D2").Value = "廖述嘉"
$ws4.Range("E2").Value = "保險期間：951811718(22年)年繳保費應繳`$71400"

$ws4.Range("A3").Value = 88
$ws4.Range("B3").Value = "中華郵政"
$ws4.Range("C3").Value = "十年快樂兒同還本終身壽險"
$ws4.Range("D3").Value = "盧秀燕"
$ws4.Range("E3").Value = "保險期間：9261025(10年）年缴保費應繳`$103431"

$ws4.Range("A4").Value = 89
$ws4.Range("B4").Value = "中華郵政"
$ws4.Range("C4").Value = "十年快樂兒同還本終身壽險"
$ws4.Range("D4").Value = "盧秀燕"
$ws4.Range("E4").Value = "保險期間：9261025(10年）年繳保費應繳`$102559"

$ws4.Range("A5").Value = 90
$ws4.Range("B5").Value = "中華郵政"
$ws4.Range("C5").Value = "金寶貝兒童保險"
$ws4.Range("D5").Value = "盧秀燕"
$ws4.Range("E5").Value = "保險期間：99824105824(6年)年繳保費應繳`$159588"

$ws4.Range("A6").Value = 91
$ws4.Range("B6").Value = "中國人壽"
$ws4.Range("C6").Value = "得意人生終身保險"
$ws4.Range("D6").Value = "盧秀燕"
$ws4.Range("E6").Value = "保險期間：9112301011230(10年)年繳保費應繳`$70309"

$ws4.Range("A7").Value = 92
$ws4.Range("B7").Value = "中國人壽"
$ws4.Range("C7").Value = "得意人生終身保險"
$ws4.Range("D7").Value = "盧秀燕"
$ws4.Range("E7").Value = "保險期間：9112301011230(10年）年繳保費應繳`$73185"

$ws4.Range("A8").Value = 93
$ws4.Range("B8").Value = "中國人壽"
$ws4.Range("C8").Value = "得意人生終身保險"
$ws4.Range("D8").Value = "盧秀燕"
$ws4.Range("E8").Value = "保險期間：9112301011230(10年)年繳保費應繳`$73234"

$ws4.Range("A9").Value = 94
$ws4.Range("B9").Value = "中國人壽"
$ws4.Range("C9").Value = "得意人生終身保險"
$ws4.Range("D9").Value = "盧秀燕"
$ws4.Range("E9").Value = "保險期間：9112301011230(10年)年缴保費應繳`$70357"

Write-Output "edit complete"
